$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.38
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("AK2").Value = 34
$ws.Range("AS2").Value = 301

# Row 7
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 6.5
$ws.Range("Y7").Value = 8.5
$ws.Range("AD7").Value = 7.5
$ws.Range("AG7").Value = 13
$ws.Range("AO7").Value = 8
$ws.Range("AP7").Value = 23
$ws.Range("AQ7").Value = 26

# Row 10
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 3.2
$ws.Range("Q10").Value = 2.15
$ws.Range("R10").Value = 1.67

# Row 11
$ws.Range("G11").Value = 3.25
$ws.Range("I11").Value = 2.15
$ws.Range("J11").Value = 3.6
$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 13
$ws.Range("X11").Value = 17
$ws.Range("Y11").Value = 11
$ws.Range("AC11").Value = 13
$ws.Range("AJ11").Value = 21
$ws.Range("AK11").Value = 17
$ws.Range("AX11").Value = 12
$ws.Range("AY11").Value = 21

# Row 12
$ws.Range("G12").Value = 3.9
$ws.Range("J12").Value = 4.75
$ws.Range("N12").Value = 7.5
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.75
$ws.Range("Q12").Value = 2.35
$ws.Range("R12").Value = 1.57
$ws.Range("U12").Value = 2.05
$ws.Range("V12").Value = 1.7
$ws.Range("X12").Value = 19
$ws.Range("AC12").Value = 7.5
$ws.Range("AM12").Value = 501
$ws.Range("AU12").Value = 9
$ws.Range("AW12").Value = 3.75
$ws.Range("AY12").Value = 26

# Row 15
$ws.Range("G15").Value = 2.1
$ws.Range("H15").Value = 2.8
$ws.Range("I15").Value = 4.1
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 1.83
$ws.Range("L15").Value = 5
$ws.Range("S15").Value = 1.67
$ws.Range("T15").Value = 2.1
$ws.Range("X15").Value = 8.5
$ws.Range("Z15").Value = 19
$ws.Range("AA15").Value = 23
$ws.Range("AB15").Value = 41
$ws.Range("AG15").Value = 8
$ws.Range("AH15").Value = 19
$ws.Range("AI15").Value = 17
$ws.Range("AN15").Value = 3.75
$ws.Range("AO15").Value = 13
$ws.Range("AT15").Value = 2.1
$ws.Range("AU15").Value = 10
$ws.Range("AX15").Value = 26

# Row 16
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 7.7
$ws.Range("I16").Value = 1.07
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 3.2
$ws.Range("L16").Value = 1.34
$ws.Range("R16").Value = 3.3
$ws.Range("S16").Value = 1.21
$ws.Range("T16").Value = 4
$ws.Range("U16").Value = 2.5
$ws.Range("V16").Value = 1.47
$ws.Range("W16").Value = 120
$ws.Range("X16").Value = 800
$ws.Range("Y16").Value = 150
$ws.Range("AA16").Value = 900
$ws.Range("AB16").Value = 450
$ws.Range("AC16").Value = 22
$ws.Range("AD16").Value = 23
$ws.Range("AE16").Value = 55
$ws.Range("AF16").Value = 250
$ws.Range("AG16").Value = 10.75
$ws.Range("AH16").Value = 6.7
$ws.Range("AI16").Value = 13.5
$ws.Range("AJ16").Value = 6
$ws.Range("AK16").Value = 11.75
$ws.Range("AL16").Value = 45
$ws.Range("AN16").Value = 26
$ws.Range("AO16").Value = 250
$ws.Range("AP16").Value = 120
$ws.Range("AT16").Value = 4
$ws.Range("AU16").Value = 12.5
$ws.Range("AV16").Value = 120
$ws.Range("AW16").Value = 3.05
$ws.Range("AX16").Value = 4.05
$ws.Range("AY16").Value = 16
$ws.Range("AZ16").Value = 7.6
$ws.Range("BA16").Value = 30
$ws.Range("BB16").Value = 250
